# Applies the "adding graphs and finetuning model" edit:
# Extends the prediction table on Sheet1 with 11 new rows (12 -> 21)
# containing day-of-week index, running idx counter, and date (column C
# formatted as a date, matching the existing style used in C2:C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data to append: day-of-week (A), idx (B), date-serial (C)
$newRows = @(
    @(6, 365, 45737),
    @(7, 366, 45738),
    @(1, 367, 45739),
    @(2, 368, 45740),
    @(3, 369, 45741),
    @(4, 370, 45742),
    @(5, 371, 45743),
    @(6, 372, 45744),
    @(7, 373, 45745),
    @(1, 374, 45746),
    @(2, 375, 45747)
)

$startRow = 11

# Reference cell whose style (date number format, numFmtId 14) the new
# column-C cells should inherit - matches C2:C10 already on the sheet.
$dateStyleSource = $ws.Cells.Item(2, 3)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $dayVal = $newRows[$i][0]
    $idxVal = $newRows[$i][1]
    $dateSerial = $newRows[$i][2]

    $ws.Cells.Item($r, 1).Value = $dayVal
    $ws.Cells.Item($r, 2).Value = $idxVal

    $cDate = $ws.Cells.Item($r, 3)
    $cDate.Value = $dateSerial

    # Copy the existing date formatting (numFmtId 14) onto the new cell
    # instead of letting a fresh NumberFormat assignment mint a new style.
    $dateStyleSource.Copy()
    $cDate.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Update the active selection to match the saved view (H17)
$ws.Range("H17").Select()
